# Establecida PilaHL como análogo de ColaHL. Actualizadas las métricas.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# --- Métricas de "Declarar Clase ListaEnlazada" (fila 6): más líneas estimadas ---
$ws.Range("B6").Value = 5

# --- Se intercambian las dos tareas de desarrollo (filas 7 y 8) ---
$ws.Range("A7").Value = "Desarrollo de metodos de la Interfaz Lista"
$ws.Range("A8").Value = "Desarrollo de comportamientos propios de Lista"

# Fila 7: "Desarrollo de metodos de la Interfaz Lista"
$ws.Range("B7").Value = 80
$ws.Range("C7").Value = 70

# Fila 8: "Desarrollo de comportamientos propios de Lista"
$ws.Range("B8").Value = 150
$ws.Range("C8").Value = 190
$ws.Range("F8").Value = 0.77430555555555547

# La fila 9 (vacía) adopta el mismo formato que las filas 7/8
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Selección activa al guardar
$ws.Range("F9").Select()
